$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.098.22'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '1.749.68'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5276'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2793'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06191'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '1.745.21'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07185'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6446'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.626'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.45'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9997'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9996'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '25.998.65'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006722'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("D21").Value = '1.968.60'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  +5.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.740'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.238'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.509'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +2.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.806'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08284'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.807'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.660'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04573'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.643'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.006'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6339'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.709'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01596'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.955'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9992'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3920'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7460'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("E45").Value = '  +3.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.348'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05352'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '31.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.77%  '
$ws.Range("E49").Value = '  +3.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.623'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3451'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.04%  '
